# Insert a new data row at row 549 (pushing the existing rows 549-586 down
# to 550-587) and populate it with the new Repollo price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 549; this shifts rows 549:586
# down to 550:587 and grows the sheet dimension from A1:R586 to A1:R587.
$ws.Rows(549).Insert()

# Populate the newly inserted row 549 with the new record's data.
$ws.Range("A549").Value = 5
$ws.Range("B549").Value = "Macroferia Regional de Talca"
$ws.Range("C549").Value = "Maule"
$ws.Range("D549").Value = 45265
$ws.Range("E549").Value = 7
$ws.Range("F549").Value = 100112006
$ws.Range("G549").Value = "Repollo"
$ws.Range("H549").Value = "Crespo record"
$ws.Range("I549").Value = "Primera"
$ws.Range("J549").Value = 4000
$ws.Range("K549").Value = 1000
$ws.Range("L549").Value = 1000
$ws.Range("M549").Value = 1000
$ws.Range("N549").Value = "$/unidad"
$ws.Range("O549").Value = "Región del Maule"
$ws.Range("P549").Value = 1000
$ws.Range("Q549").Value = 1
$ws.Range("R549").Value = "Hortaliza"
